# Prignano.xlsx update: append daily COVID tracking rows through 2021-12-08 (8/12)
# (commit message: "aggiornamento fino a 8/12")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New data rows 386..464: columns B (nuovi pos.), C (somma mobile 7gg.),
# D (somma mobile 7gg. per 100mila abitanti). Column A is the date serial,
# which continues sequentially (44460 .. 44538) from the prior last row (44459).
$newData = @(
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(2, 2, 53.53319057815846),
    @(0, 2, 53.53319057815846),
    @(0, 2, 53.53319057815846),
    @(0, 2, 53.53319057815846),
    @(0, 2, 53.53319057815846),
    @(0, 2, 53.53319057815846),
    @(1, 3, 80.29978586723769),
    @(0, 1, 26.76659528907923),
    @(0, 1, 26.76659528907923),
    @(0, 1, 26.76659528907923),
    @(0, 1, 26.76659528907923),
    @(0, 1, 26.76659528907923),
    @(0, 1, 26.76659528907923),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(0, 0, 0),
    @(5, 5, 133.8329764453961),
    @(0, 5, 133.8329764453961),
    @(0, 5, 133.8329764453961),
    @(0, 5, 133.8329764453961),
    @(0, 5, 133.8329764453961),
    @(0, 5, 133.8329764453961),
    @(0, 5, 133.8329764453961),
    @(0, 0, 0),
    @(6, 6, 160.5995717344754),
    @(0, 6, 160.5995717344754),
    @(0, 6, 160.5995717344754),
    @(0, 6, 160.5995717344754),
    @(0, 6, 160.5995717344754),
    @(1, 7, 187.3661670235546),
    @(0, 7, 187.3661670235546),
    @(0, 1, 26.76659528907923),
    @(5, 6, 160.5995717344754),
    @(2, 8, 214.1327623126338),
    @(1, 9, 240.8993576017131),
    @(0, 9, 240.8993576017131),
    @(2, 10, 267.6659528907923),
    @(1, 11, 294.4325481798715),
    @(0, 11, 294.4325481798715)
)

$firstNewRow = 386
$lastOldRow = $firstNewRow - 1
$startDate = $ws.Cells.Item($lastOldRow, 1).Value2 + 1

for ($i = 0; $i -lt $newData.Count; $i++) {
    $r = $firstNewRow + $i
    $ws.Cells.Item($r, 1).Value2 = $startDate + $i
    $ws.Cells.Item($r, 2).Value2 = $newData[$i][0]
    $ws.Cells.Item($r, 3).Value2 = $newData[$i][1]
    $ws.Cells.Item($r, 4).Value2 = $newData[$i][2]
}

# Carry the date-column formatting (bold, boxed border, center/top aligned,
# yyyy-mm-dd date number format) down from the previous last row, matching the
# style already applied to A2:A385.
$lastNewRow = $firstNewRow + $newData.Count - 1
$srcCell = $ws.Cells.Item($lastOldRow, 1)
$destRange = $ws.Range($ws.Cells.Item($firstNewRow, 1), $ws.Cells.Item($lastNewRow, 1))
$srcCell.Copy()
$destRange.PasteSpecial(-4122)
$excel.CutCopyMode = $false
